$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "626.38") need to be
# forced to Text format first, otherwise Excel auto-converts the assigned
# string into a numeric value instead of keeping it as text (as in the source).
$textForceCells = @(
    'D5', 'D6', 'D8', 'D11', 'D12', 'D14', 'D21', 'D23', 'D25', 'D28', 'D31', 'D33', 'D34', 'D40', 'D41', 'D42', 'D46', 'D47', 'D48'
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply every updated cell value from the refreshed crypto price feed.
$ws.Range('D2').Value = '69.228.18'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').Value = '3.775.05'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '626.38'
$ws.Range('E5').Value = '  +4.39%  '
$ws.Range('D6').Value = '166.09'
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('D7').Value = '3.773.80'
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('E10').Value = '  +2.16%  '
$ws.Range('D11').Value = '0.461'
$ws.Range('E11').Value = '  +3.60%  '
$ws.Range('D12').Value = '6.75'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').Value = '35.86'
$ws.Range('E14').Value = '  +2.19%  '
$ws.Range('D15').Value = '4.412.52'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '3.774.88'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '69.232.90'
$ws.Range('E17').Value = '  +2.16%  '
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('D21').Value = '468.47'
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('D23').Value = '0.708'
$ws.Range('E23').Value = '  +2.44%  '
$ws.Range('E24').Value = '  +4.30%  '
$ws.Range('D25').Value = '83.22'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  +2.07%  '
$ws.Range('E27').Value = '  +4.25%  '
$ws.Range('D28').Value = '10.07'
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '3.925.00'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').Value = '2.68'
$ws.Range('E31').Value = '  +3.50%  '
$ws.Range('E32').Value = '  +2.61%  '
$ws.Range('D33').Value = '7.19'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').Value = '28.83'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  +15.09%  '
$ws.Range('D37').Value = '3.727.46'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('E39').Value = '  +2.65%  '
$ws.Range('D40').Value = '3.39'
$ws.Range('E40').Value = '  +7.50%  '
$ws.Range('D41').Value = '5.83'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').Value = '0.969'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').Value = '43.13'
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '152.92'
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('D48').Value = '46.73'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('E49').Value = '  +4.43%  '
$ws.Range('E50').Value = '  +1.75%  '
$ws.Range('E51').Value = '  -0.20%  '

# Drop the Text format back to the default cell style now that the literal
# string is committed, so formatting matches the rest of the sheet.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
